$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.049.10'
$ws.Range("E2").Value = '  -2.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.993.14'
$ws.Range("E3").Value = '  -1.14%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.014'
$ws.Range("E4").Value = '  +0.24%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '329.97'
$ws.Range("E5").Value = '  -0.65%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.013'
$ws.Range("E6").Value = '  +0.24%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4957'
$ws.Range("E7").Value = '  -1.75%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4177'
$ws.Range("E8").Value = '  -2.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '55.03'
$ws.Range("E9").Value = '  +1.39%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08862'
$ws.Range("E10").Value = '  -4.09%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.089'
$ws.Range("E11").Value = '  -3.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '22.96'
$ws.Range("E12").Value = '  -2.92%  '

$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.971'
$ws.Range("E13").Value = '  -2.07%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.951.65'
$ws.Range("E14").Value = '  -0.78%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.399'
$ws.Range("E15").Value = '  -2.55%  '

$ws.Range("E16").Value = '  +0.30%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '92.22'
$ws.Range("E17").Value = '  -4.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001102'
$ws.Range("E18").Value = '  -2.36%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06671'
$ws.Range("E19").Value = '  -0.09%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.40'
$ws.Range("E20").Value = '  -2.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.013'
$ws.Range("E21").Value = '  +0.45%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.959'
$ws.Range("E22").Value = '  -0.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '29.083.20'
$ws.Range("E23").Value = '  -2.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.95'
$ws.Range("E24").Value = '  -1.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.326'
$ws.Range("E25").Value = '  +2.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.237.79'
$ws.Range("E26").Value = '  +1.24%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.78'
$ws.Range("E27").Value = '  -0.37%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '156.98'
$ws.Range("E28").Value = '  -1.77%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.239'
$ws.Range("E29").Value = '  -3.58%  '

$ws.Range("E30").Value = '  -4.87%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '126.79'
$ws.Range("E31").Value = '  -1.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.041'
$ws.Range("E32").Value = '  -2.21%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09873'
$ws.Range("E33").Value = '  -1.22%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.525'
$ws.Range("E34").Value = '  -4.44%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.826'
$ws.Range("E35").Value = '  -1.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.747'
$ws.Range("E36").Value = '  -1.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02409'

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.308'
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '9.045'
$ws.Range("E39").Value = '  -6.33%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06377'
$ws.Range("E40").Value = '  -0.56%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6459'
$ws.Range("E41").Value = '  -2.18%  '

$ws.Range("E42").Value = '  -2.63%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1970'
$ws.Range("E43").Value = '  -5.64%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.013'

$ws.Range("B45").Value = 'WEMIXTOKEN'
$ws.Range("C45").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.357'
$ws.Range("E45").Value = '  +5.02%  '

$ws.Range("B46").Value = 'Decentraland'
$ws.Range("C46").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6165'
$ws.Range("E46").Value = '  -3.47%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '13.27'
$ws.Range("E47").Value = '  -2.79%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.161'
$ws.Range("E48").Value = '  -2.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00000000348'
$ws.Range("E49").Value = '  +8.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.494'
$ws.Range("E50").Value = '  -1.17%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.161'
$ws.Range("E51").Value = '  +5.98%  '

